$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember label cells before editing so new rows can reuse the same shared strings.
$labelA3 = $ws.Range("A3").Value()
$labelA4 = $ws.Range("A4").Value()
$labelA5 = $ws.Range("A5").Value()

# --- Row 2: update B2 and D2 ---
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 10

# --- New rows 7, 8, 9 (values entered first so new shared strings land in this order) ---
$ws.Range("A9").Value = $labelA3
$ws.Range("B9").Value = "Повысить обороты"

# --- Row 3: change result text ---
$ws.Range("B3").Value = "Снизить напряжение"

# --- Row 6: change result text ---
$ws.Range("B6").Value = "Обратиться к специалисту"

# --- Row 4: keep only A4/B4, clear C4/D4 ---
$ws.Range("B4").Value = 37
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

# --- Row 5: keep only A5/B5, clear C5/D5 ---
$ws.Range("B5").Value = 10
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()

$ws.Range("A7").Value = $labelA4
$ws.Range("B7").Value = 36

$ws.Range("A8").Value = $labelA5
$ws.Range("B8").Value = 4

# --- Selection matches the final saved state ---
$ws.Range("B9").Select()
